# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on the Overview
#    sheet (E2, F2) and on each language sheet's Status cell (C2).
# 2) Narrow the "Latest Handoff/Handback Datetime" style columns from
#    ~17.22 chars to ~13.41 chars: Overview!E:F and the C column on the
#    zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1) Update status text -------------------------------------------------

$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# --- 2) Narrow the datetime columns ----------------------------------------
# The workbook's column widths were authored outside Excel, so the raw
# OOXML "width" values (e.g. 17.2159881591797) don't align to Excel's
# internal pixel grid (width is always a multiple of 1/6 once it has been
# round-tripped through the ColumnWidth property). 12.5 is the
# ColumnWidth ("characters") value that lands on the pixel bucket closest
# to the target width of 13.4101848602295 characters.

$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F
$zhcn.Columns.Item(3).ColumnWidth     = $newColumnWidth  # column C
$dede.Columns.Item(3).ColumnWidth     = $newColumnWidth  # column C
